$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new game data (challenge replaced)
$ws.Range("A2").Value = "bvarNUA0"
$ws.Range("E2").Value = "https://lichess.org/bvarNUA0"
$ws.Range("F2").Value = 4087

# Remove rows 3 and 4 (clear their contents, keep empty rows)
$ws.Range("A3:H4").ClearContents()
